$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("M2").Value = 10.23707633333333
$ws.Range("N2").Value = 30.711229
$ws.Range("O2").Value = 0.2523196022162781
$ws.Range("P2").Value = 0.2523196022162781
$ws.Range("Q2").Value = 3.065949764092889
$ws.Range("R2").Value = 27.593547876836
$ws.Range("S2").Value = 0.1258936472501387
$ws.Range("T2").Value = 0.1258936472501387
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("O3").Value = 0.5485767480449855
$ws.Range("P3").Value = 0.5485767480449855
$ws.Range("Q3").Value = 6.665787106836445
$ws.Range("S3").Value = 0.2737097197418937
$ws.Range("T3").Value = 0.2737097197418937
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("O4").Value = 0.1991036497387364
$ws.Range("P4").Value = 0.1991036497387364
$ws.Range("S4").Value = 0.09934180470425039
$ws.Range("T4").Value = 0.0993418047042504
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("M5").Value = 10.23707633333333
$ws.Range("N5").Value = 30.711229
$ws.Range("O5").Value = 0.2523196022162781
$ws.Range("P5").Value = 0.2523196022162781
$ws.Range("Q5").Value = 3.078913315089666
$ws.Range("R5").Value = 27.710219835807
$ws.Range("S5").Value = 0.1264259549661395
$ws.Range("T5").Value = 0.1264259549661395
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("O6").Value = 0.5485767480449855
$ws.Range("P6").Value = 0.5485767480449855
$ws.Range("S6").Value = 0.2748670283030917
$ws.Range("T6").Value = 0.2748670283030917
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("O7").Value = 0.1991036497387364
$ws.Range("P7").Value = 0.1991036497387364
$ws.Range("S7").Value = 0.09976184503448603
$ws.Range("T7").Value = 0.09976184503448603
